# "draw(layout): add timetable, attendance"
#
# Adds a new "courses" mini table under column I (rows 21-24: courses / id /
# subject_id / weekday), using a new orange fill style for the table header,
# and clears the stray "subject_id" value that used to sit at I17 (the old
# end of the "schedules" table) now that "subject_id" has moved into the new
# "courses" block. Also nudges the sheet's zoom level/selection and touches
# column K's width, matching the reviewed layout pass.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "subject_id" label that used to trail the schedules table (I17) is
# removed - it now lives in the new courses table instead.
$ws.Range("I17").ClearContents()

# New "courses" table header (orange fill) + its three fields.
$ws.Range("I21").Value = "courses"
$ws.Range("I21").Interior.Color = 0x0000C0FF

$ws.Range("I22").Value = "id"
$ws.Range("I23").Value = "subject_id"
$ws.Range("I24").Value = "weekday"

# Column K (attendances table) picks up an explicit width while reviewing
# the new layout.
$ws.Columns("K").ColumnWidth = 11.59

# Zoom out a bit and leave the selection on the new courses table.
$excel.ActiveWindow.Zoom = 115
$ws.Range("L21").Select() | Out-Null
